$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(99).Insert()

$ws.Cells.Item(99, 1).Value = 3
$ws.Cells.Item(99, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(99, 3).Value = "Coquimbo"
$ws.Cells.Item(99, 4).Value = 44588
$ws.Cells.Item(99, 5).Value = 5
$ws.Cells.Item(99, 6).Value = 100112030
$ws.Cells.Item(99, 7).Value = "Poroto granado"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 78
$ws.Cells.Item(99, 11).Value = 25000
$ws.Cells.Item(99, 12).Value = 26000
$ws.Cells.Item(99, 13).Value = 25487
$ws.Cells.Item(99, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 1019
$ws.Cells.Item(99, 17).Value = 25
$ws.Cells.Item(99, 18).Value = "Hortaliza"
